# Insert a new data row at spreadsheet row 65 (pushing the existing
# rows 65..147 down to 66..148) and populate the new row with the
# latest weekly price observation for this market/variety.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(65).Insert()

$ws.Cells.Item(65, 1).Value  = 8
$ws.Cells.Item(65, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(65, 3).Value  = "Coquimbo"
$ws.Cells.Item(65, 4).Value  = [datetime]"2021-09-28"
$ws.Cells.Item(65, 5).Value  = 4
$ws.Cells.Item(65, 6).Value  = 100112003
$ws.Cells.Item(65, 7).Value  = "Ajo"
$ws.Cells.Item(65, 8).Value  = "Chino"
$ws.Cells.Item(65, 9).Value  = "Primera"
$ws.Cells.Item(65, 10).Value = 520
$ws.Cells.Item(65, 11).Value = 16000
$ws.Cells.Item(65, 12).Value = 17000
$ws.Cells.Item(65, 13).Value = 16500
$ws.Cells.Item(65, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(65, 15).Value = "China"
$ws.Cells.Item(65, 16).Value = 1650
$ws.Cells.Item(65, 17).Value = 10
$ws.Cells.Item(65, 18).Value = "Hortaliza"
